# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 138
$ws1.Range("F4").Value = 432
$ws1.Range("F5").Value = 1739
$ws1.Range("F11").Value = 4947
$ws1.Range("F13").Value = 38
$ws1.Range("F17").Value = 187
$ws1.Range("F21").Value = 3904
$ws1.Range("F23").Value = 676
$ws1.Range("F28").Value = 24
$ws1.Range("F34").Value = 965
$ws1.Range("F35").Value = 2474

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 138
$ws4.Range("F4").Value = 432
$ws4.Range("F5").Value = 1739
$ws4.Range("F11").Value = 4947
$ws4.Range("F13").Value = 38
$ws4.Range("F17").Value = 187
$ws4.Range("F21").Value = 3904
$ws4.Range("F23").Value = 676
$ws4.Range("F28").Value = 24
$ws4.Range("F35").Value = 965
$ws4.Range("F36").Value = 2474
